# Used Addons.docx edit script
# -resolved #30 -resolved #58 -resolved #206 -resolved #212 -resolved #213
#
# Summary of changes:
#  * Neutral/General table: remove "Standard Animations / Sanctuary / GRAA" row
#  * USA Vehicles table: remove "A-10 / Diesel" and "M60A1 / Inquisitor" rows
#  * USA Weapons table: remove "M1911 / SJB", "M67 Frag Grenade, AN-M8, M18 / KyleSarnik / BD Grenade Pack",
#       "M60 / Laser", "M47 Dragon / Combat (CBT)", "PRC-119 / Solus / SLX", "M21 / SJB" rows
#  * USSR Vehicles table: remove "Mi-8MTV / OUR WEAPONS Project / OWP Mi8 High Quality",
#       "ZSU-23-4 / ORCS", "T-54 / RedHammerStudios" rows
#  * USSR Weapons table: rename "AK-74" row to "AKM GP-25" (author stays RedHammerStudios),
#       remove "RGO Hand Grenade", "PKM", "RPG-7", "RPG-18", "RDG1 Hand Grenade",
#       "AKM", old "AKM GP-25", "RPG-2" rows
#  * FIA section: remove the whole "Vehicles" heading paragraph and its table
#       (ZIL-130/131 / Our Weapons Project)
#
# NOTE: table / range object handles returned by the COM layer go stale as soon as
# *any* part of the document mutates, so every helper below re-resolves the table
# (and row/cell) fresh from $d.Tables right before it touches it.

$d = $word.ActiveDocument

function Remove-TableRowByFirstCell($tableIndex, $text) {
    $marker = $text + "`r`a"
    $rowCount = $d.Tables.Item($tableIndex).Rows.Count
    for ($r = $rowCount; $r -ge 1; $r--) {
        $cellText = $d.Tables.Item($tableIndex).Cell($r, 1).Range.Text
        if ($cellText -eq $marker) {
            $d.Tables.Item($tableIndex).Rows.Item($r).Delete()
            return $true
        }
    }
    return $false
}

# --- Rename AK-74 -> AKM GP-25 in the USSR Weapons table (table 7) ---
$d.Tables.Item(7).Cell(2, 1).Range.Find.Execute("AK-74", $true, $false, $false, $false, $false, `
    $true, 1, $false, "AKM GP-25", 2) | Out-Null

# --- Neutral/General table (table 1) ---
Remove-TableRowByFirstCell 1 "Standard Animations" | Out-Null

# --- USA Vehicles table (table 3) ---
Remove-TableRowByFirstCell 3 "M60A1" | Out-Null
Remove-TableRowByFirstCell 3 "A-10" | Out-Null

# --- USA Weapons table (table 4) ---
Remove-TableRowByFirstCell 4 "M21" | Out-Null
Remove-TableRowByFirstCell 4 "PRC-119" | Out-Null
Remove-TableRowByFirstCell 4 "M47 Dragon" | Out-Null
Remove-TableRowByFirstCell 4 "M60" | Out-Null
Remove-TableRowByFirstCell 4 "M67 Frag Grenade, AN-M8, M18" | Out-Null
Remove-TableRowByFirstCell 4 "M1911" | Out-Null

# --- USSR Vehicles table (table 6) ---
Remove-TableRowByFirstCell 6 "T-54" | Out-Null
Remove-TableRowByFirstCell 6 "ZSU-23-4" | Out-Null
Remove-TableRowByFirstCell 6 "Mi-8MTV" | Out-Null

# --- USSR Weapons table (table 7), remove obsolete rows ---
Remove-TableRowByFirstCell 7 "RPG-2" | Out-Null
Remove-TableRowByFirstCell 7 "AKM GP-25" | Out-Null
Remove-TableRowByFirstCell 7 "AKM" | Out-Null
Remove-TableRowByFirstCell 7 "RDG1 Hand Grenade" | Out-Null
Remove-TableRowByFirstCell 7 "RPG-18" | Out-Null
Remove-TableRowByFirstCell 7 "RPG-7" | Out-Null
Remove-TableRowByFirstCell 7 "PKM" | Out-Null
Remove-TableRowByFirstCell 7 "RGO Hand Grenade" | Out-Null

# --- FIA Vehicles section: remove the table and its "Vehicles" heading paragraph ---
$d.Tables.Item(9).Delete()

$headingIndex = 0
$paraCount = $d.Content.Paragraphs.Count
for ($i = 1; $i -le $paraCount; $i++) {
    if ($d.Content.Paragraphs.Item($i).Range.Text -eq "Vehicles`r") {
        $headingIndex = $i
    }
}
if ($headingIndex -gt 0) {
    $d.Content.Paragraphs.Item($headingIndex).Range.Delete()
}
